# Armor.xlsx — "Added some new entries for our Armor.csv file"
# Appends 15 new armor/outfit rows (rows 9-23) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Name (col A), weight (col B), item hp (col C),
# baseid (col D), Value (col E), damage rating (col F)
$newArmor = @(
    @{ Row = 9;  Name = "Vault 77 Jumpsuit";           Weight = 6;   Hp = 1;  BaseId = "000cafbe"; Value = 100;  Dmg = 1 },
    @{ Row = 10; Name = "Vault 87 Jumpsuit";           Weight = 6;   Hp = 1;  BaseId = "000340ed"; Value = 100;  Dmg = 1 },
    @{ Row = 11; Name = "Vault 92 Jumpsuit";           Weight = 6;   Hp = 1;  BaseId = "000b73f3"; Value = 100;  Dmg = 1 },
    @{ Row = 12; Name = "Vault 101 Jumpsuit";          Weight = 6;   Hp = 1;  BaseId = "0000431e"; Value = 100;  Dmg = 1 },
    @{ Row = 13; Name = "Vault 106 Jumpsuit";          Weight = 6;   Hp = 1;  BaseId = "000b73f2"; Value = 100;  Dmg = 1 },
    @{ Row = 14; Name = "Vault 108 Jumpsuit";          Weight = 6;   Hp = 1;  BaseId = "000b73f1"; Value = 100;  Dmg = 1 },
    @{ Row = 15; Name = "Vault 112 Jumpsuit";          Weight = 6;   Hp = 1;  BaseId = "000340ef"; Value = 100;  Dmg = 1 },
    @{ Row = 16; Name = "Vault 101 Utility Jumpsuit";  Weight = 8;   Hp = 2;  BaseId = "000425ba"; Value = 100;  Dmg = 1 },
    @{ Row = 17; Name = "Vault 101 Security Armor";    Weight = 70;  Hp = 15; BaseId = "0003411c"; Value = 100;  Dmg = 12 },
    @{ Row = 18; Name = "Vault Lab Uniform";           Weight = 6;   Hp = 1;  BaseId = "0001cbdc"; Value = 100;  Dmg = 2 },
    @{ Row = 19; Name = "Combat Armor";                Weight = 390; Hp = 25; BaseId = "00020420"; Value = 400;  Dmg = 32 },
    @{ Row = 20; Name = "Ranger Battle Armor";         Weight = 430; Hp = 27; BaseId = "00023030"; Value = 1100; Dmg = 39 },
    @{ Row = 21; Name = "Rivet City security uniform"; Weight = 330; Hp = 20; BaseId = "000239cc"; Value = 100;  Dmg = 24 },
    @{ Row = 22; Name = "Talon Combat Armor";          Weight = 275; Hp = 25; BaseId = "000a6f76"; Value = 300;  Dmg = 28 },
    @{ Row = 23; Name = "Tennypenny Security Uniform"; Weight = 180; Hp = 20; BaseId = "00034119"; Value = 100;  Dmg = 24 }
)
$byRow = @{}
foreach ($entry in $newArmor) { $byRow[$entry.Row] = $entry }

# The text columns (Name / baseid) were pasted in from the source CSV in a few
# uneven batches rather than strictly row-by-row, so replay that exact
# sequence of (row, column) writes to land each string in the same slot of
# the shared-string table as the saved workbook.
$textWriteOrder = @(
    @(9,'A'), @(9,'D'),
    @(10,'A'), @(11,'A'), @(12,'A'), @(13,'A'), @(14,'A'), @(15,'A'),
    @(10,'D'), @(11,'D'), @(12,'D'), @(13,'D'), @(14,'D'), @(15,'D'),
    @(16,'A'), @(16,'D'),
    @(17,'A'), @(17,'D'),
    @(18,'A'), @(18,'D'),
    @(19,'A'), @(19,'D'),
    @(20,'D'), @(20,'A'),
    @(21,'D'), @(21,'A'),
    @(22,'A'), @(22,'D'),
    @(23,'D'), @(23,'A')
)

foreach ($pair in $textWriteOrder) {
    $row = $pair[0]
    $col = $pair[1]
    $entry = $byRow[$row]
    $value = if ($col -eq 'A') { $entry.Name } else { $entry.BaseId }
    $ws.Range("$col$row").Value = $value
}

# Numeric columns don't touch the shared-string table, so plain row order is fine.
foreach ($entry in $newArmor) {
    $ws.Range("B$($entry.Row)").Value = $entry.Weight
    $ws.Range("C$($entry.Row)").Value = $entry.Hp
    $ws.Range("E$($entry.Row)").Value = $entry.Value
    $ws.Range("F$($entry.Row)").Value = $entry.Dmg
}

# Column A widens slightly to fit the longer new names (bestFit recalculation)
$ws.Columns.Item(1).ColumnWidth = 26.92

# Leave the same cell selected as in the saved workbook
[void]$ws.Range("E16").Select()
